# New weekly price observation for "Locoto" (Vega Modelo de Temuco) is
# inserted at the top of the price-history block (row 87), pushing the
# existing rows 87-105 down to 88-106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 87, shifting rows 87:105 down to 88:106.
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new observation.
$ws.Range("A87").Value = 10
$ws.Range("B87").Value = "Vega Modelo de Temuco"
$ws.Range("C87").Value = "La Araucanía"
$ws.Range("D87").Value = 45204
$ws.Range("E87").Value = 9
$ws.Range("F87").Value = 100112042
$ws.Range("G87").Value = "Locoto"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 100
$ws.Range("K87").Value = 2200
$ws.Range("L87").Value = 2200
$ws.Range("M87").Value = 2200
$ws.Range("N87").Value = "$/kilo"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 2200
$ws.Range("Q87").Value = 1
$ws.Range("R87").Value = "Hortaliza"
